{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldText = \"Cerradura configureCerradura(int numEstructuras); // Modificado\";\nconst newText = \"Cerradura configureCerradura(const std::vector<int>& tamanos); // Modificado\";\n\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === oldText) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not find target paragraph: \" + oldText);\n}\n\nconst targetParagraph = paragraphs.items[targetIndex];\n\n// Update the text of the run/paragraph in-place (keeps run formatting).\ntargetParagraph.insertText(newText, \"Replace\");\n\n// The following empty paragraph is removed (its paragraph mark is merged\n// into the edited paragraph above), matching the diff which drops that\n// whole <w:p> block.\nif (targetIndex + 1 < paragraphs.items.length) {\n  const nextParagraph = paragraphs.items[targetIndex + 1];\n  nextParagraph.load(\"text\");\n  await context.sync();\n  if (nextParagraph.text === \"\") {\n    nextParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"Cerradura configureCerradura(int numEstructuras); // Modificado\"\n$newText = \"Cerradura configureCerradura(const std::vector<int>& tamanos); // Modificado\"\n\n# Locate the paragraph containing the old signature text.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq $oldText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find target paragraph with text: $oldText\"\n}\n\n# Replace the run text, keeping the paragraph mark (so formatting/rPr persists).\n$target.Range.Text = $newText\n\n# Remove the following empty paragraph entirely (merges it away), matching\n# the diff which deletes that whole <w:p> block.\n$nextPara = $target.Next()\nif ($nextPara -ne $null -and $nextPara.Range.Text.TrimEnd(\"`r`a\") -eq \"\") {\n    $nextPara.Range.Delete() | Out-Null\n}\n"}
